$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1. Table on slide 6: switch the applied table style from the
#    custom "Table_0" style to the built-in style
#    {A9718B63-6C17-444E-B8D4-17D96AA3EFB3}.
# ------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A9718B63-6C17-444E-B8D4-17D96AA3EFB3}")
    }
}

# ------------------------------------------------------------------
# 2. Presentation theme: swap the "Integral" colour scheme that is
#    currently applied to the slide master/deck for the default
#    "Office Theme" colour scheme (the two theme parts effectively
#    trade their colour palettes).
# ------------------------------------------------------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72
